# Regenerate the s_vals data (filtering save games) for this player's sheet.
# Each row's B:E columns hold recomputed stat ratios and G holds their sum.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @{ B = 3.230985683306322;     C = 0.3127903958511391;    D = 0.8054896365839992;  E = 8.660232485948974 }
    3  = @{ B = 0.01514828764759746;   C = 3099.503889238888;     D = 26.21740644021617;   E = 645.3272768299601 }
    4  = @{ B = 3.230985683306322;     C = 1.667794583268128;     D = 0.1575252929769615;  E = 0.496779210170732 }
    5  = @{ B = 0.0000002070225997297115; C = 0.000000000002643885110842348; D = 26.21740644021617;   E = 645.3272768299601 }
    6  = @{ B = 0.01514828764759746;   C = 0.002777888934908601;  D = 0.8054896365839992;  E = 0.496779210170732 }
    7  = @{ B = 3.230985683306322;     C = 1.667794583268128;     D = 0.8054896365839992;  E = 0.496779210170732 }
    8  = @{ B = 0.3048080303191223;    C = 1.667794583268128;     D = 0.8054896365839992;  E = 645.3272768299601 }
    9  = @{ B = 3.230985683306322;     C = 1.667794583268128;     D = 0.1575252929769615;  E = 0.496779210170732 }
    10 = @{ B = 0.003994804209775715;  C = 1.667794583268128;     D = 49627605961.23487;    E = 2367095152636972 }
    11 = @{ B = 0.00002074986032285508; C = 0.00007097389502863649; D = 337.1190423067083;   E = 645.3272768299601 }
}

$G = @{
    2  = 13.00949820169043
    3  = 3771.063720796712
    4  = 5.553084769722144
    5  = 671.5446834772015
    6  = 1.320195023337237
    7  = 6.201049113329182
    8  = 648.1053690801313
    9  = 5.553084769722144
    10 = 2367144780242934
    11 = 982.4464108604236
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals.B
    $ws.Cells.Item($row, 3).Value = $vals.C
    $ws.Cells.Item($row, 4).Value = $vals.D
    $ws.Cells.Item($row, 5).Value = $vals.E
    $ws.Cells.Item($row, 7).Value = $G[$row]
}
